$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-12-24 Wednesday" "2025-12-25 Thursday"

Replace-Text "322×4=1288" "109×4=436"
Replace-Text "997×5=4985" "740×9=6660"
Replace-Text "590×9=5310" "881×6=5286"
Replace-Text "544×5=2720" "204×4=816"
Replace-Text "529×9=4761" "905×6=5430"
Replace-Text "952×8=7616" "519×7=3633"
Replace-Text "403×7=2821" "135×3=405"
Replace-Text "955×4=3820" "639×7=4473"
Replace-Text "102×7=714" "382×6=2292"
Replace-Text "909×3=2727" "871×8=6968"
Replace-Text "606×6=3636" "122×3=366"
Replace-Text "168×8=1344" "431×7=3017"
Replace-Text "654×3=1962" "586×7=4102"
Replace-Text "613×7=4291" "158×3=474"
Replace-Text "826×5=4130" "607×7=4249"
Replace-Text "744×9=6696" "152×8=1216"
Replace-Text "452×3=1356" "465×6=2790"
Replace-Text "212×6=1272" "570×9=5130"
Replace-Text "108×7=756" "598×5=2990"
Replace-Text "162×9=1458" "582×7=4074"
Replace-Text "723×4=2892" "154×4=616"
Replace-Text "695×7=4865" "752×8=6016"
Replace-Text "806×6=4836" "981×7=6867"
Replace-Text "659×2=1318" "276×8=2208"
Replace-Text "530×2=1060" "818×2=1636"
